# Added some scenarios for class
#
# 1. Rename the "ClassDetailsForm" sheet to "Class"
# 2. Change the ClassTopic value on row 2 (C2) from "Calculus3333" to "AI_"
# 3. Add two new test-data rows (11 and 12) for new scenarios:
#      - InvalidDataForMandatoryFieldsForEdit
#      - SpecialCharValidationForEdit
# 4. Widen column A to fit the new (longer) TestCase names
# 5. Leave the final selection on D19, matching the author's last cursor spot

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClassDetailsForm")

# --- 1. Rename sheet ---
$ws.Name = "Class"

# --- 2. Update C2 (ClassTopic for the first scenario) ---
$ws.Range("C2").Value = "AI_"

# --- 3a. Row 11: InvalidDataForMandatoryFieldsForEdit scenario ---
$ws.Range("A11").Style = $ws.Range("A10").Style
$ws.Range("A11").Value = "InvalidDataForMandatoryFieldsForEdit"

$ws.Range("D11").Value = 234

$ws.Range("E11").Style = $ws.Range("E2").Style
$ws.Range("E11").Value = "03/20/2025,03/21/2025"

$ws.Range("G11").Style = $ws.Range("G2").Style
$ws.Range("G11").Value = "Saran"

$ws.Range("H11").Style = $ws.Range("H3").Style
$ws.Range("H11").Value = "Active"

$ws.Range("I11").Style = $ws.Range("I4").Style
$ws.Range("I11").Value = "asd"

$ws.Range("J11").Style = $ws.Range("J4").Style
$ws.Range("J11").Value = "no"

$ws.Range("K11").Style = $ws.Range("K4").Style
$ws.Range("K11").Value = "no"

$ws.Range("L11").Value = "Error"

$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(10).RowHeight

# --- 3b. Row 12: SpecialCharValidationForEdit scenario ---
$ws.Range("A12").Style = $ws.Range("A10").Style
$ws.Range("A12").Value = "SpecialCharValidationForEdit"

$ws.Range("D12").Value = "%$*"

$ws.Range("E12").Style = $ws.Range("E2").Style
$ws.Range("E12").Value = "03/20/2025,03/21/2025"

$ws.Range("G12").Style = $ws.Range("G2").Style
$ws.Range("G12").Value = "#$*"

$ws.Range("H12").Style = $ws.Range("H3").Style
$ws.Range("H12").Value = "Active"

$ws.Range("I12").Style = $ws.Range("I4").Style
$ws.Range("I12").Value = "#*^"

$ws.Range("L12").Value = "Error"

$ws.Rows.Item(12).RowHeight = $ws.Rows.Item(10).RowHeight

# --- 4. Widen column A so the longer scenario names fit ---
$ws.Columns.Item(1).ColumnWidth = 32.15

# --- 5. Restore last-used selection ---
$ws.Range("D19").Select()
